$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1875
$ws.Range("C2").Value = 0.5625
$ws.Range("J2").Value = 0.0625
$ws.Range("P2").Value = 0.125
$ws.Range("S2").Value = 0.0625
$ws.Range("B3").Value = 0.1
$ws.Range("C3").Value = 0.1
$ws.Range("P3").Value = 0.6
$ws.Range("S3").Value = 0.2
$ws.Range("J4").Value = 0.1666666666666667
$ws.Range("P4").Value = 0.8333333333333334
$ws.Range("J6").Value = 0.6666666666666666
$ws.Range("Q6").Value = 0.1666666666666667
$ws.Range("S6").Value = 0.1666666666666667
$ws.Range("B7").Value = 0.1428571428571428
$ws.Range("D7").Value = 0.07142857142857142
$ws.Range("J7").Value = 0.2142857142857143
$ws.Range("O7").Value = 0.1428571428571428
$ws.Range("Q7").Value = 0.1428571428571428
$ws.Range("S7").Value = 0.2857142857142857
$ws.Range("B8").Value = 0.1333333333333333
$ws.Range("D8").Value = 0.03333333333333333
$ws.Range("E8").Value = 0.03333333333333333
$ws.Range("J8").Value = 0.1666666666666667
$ws.Range("O8").Value = 0.03333333333333333
$ws.Range("Q8").Value = 0.2
$ws.Range("R8").Value = 0.1666666666666667
$ws.Range("S8").Value = 0.2333333333333333
$ws.Range("D9").Value = 0.125
$ws.Range("Q9").Value = 0.375
$ws.Range("R9").Value = 0.125
$ws.Range("S9").Value = 0.375
$ws.Range("B10").Value = 0.05454545454545454
$ws.Range("D10").Value = 0.02727272727272727
$ws.Range("F10").Value = 0.01818181818181818
$ws.Range("J10").Value = 0.1181818181818182
$ws.Range("O10").Value = 0.01818181818181818
$ws.Range("Q10").Value = 0.3454545454545455
$ws.Range("R10").Value = 0.1
$ws.Range("S10").Value = 0.3181818181818182
$ws.Range("G11").Value = 0.1071428571428571
$ws.Range("J11").Value = 0.1071428571428571
$ws.Range("K11").Value = 0.1785714285714286
$ws.Range("L11").Value = 0.6071428571428571
$ws.Range("G12").Value = 0.6111111111111112
$ws.Range("J12").Value = 0.3333333333333333
$ws.Range("K12").Value = 0.05555555555555555
$ws.Range("J13").Value = 1
$ws.Range("H15").Value = 0.1333333333333333
$ws.Range("J15").Value = 0.4
$ws.Range("K15").Value = 0.06666666666666667
$ws.Range("O15").Value = 0.06666666666666667
$ws.Range("S15").Value = 0.3333333333333333
$ws.Range("H16").Value = 0.1538461538461539
$ws.Range("I16").Value = 0.07692307692307693
$ws.Range("J16").Value = 0.3846153846153846
$ws.Range("O16").Value = 0.07692307692307693
$ws.Range("S16").Value = 0.3076923076923077
$ws.Range("F17").Value = 0.04
$ws.Range("H17").Value = 0.2
$ws.Range("I17").Value = 0.02
$ws.Range("J17").Value = 0.36
$ws.Range("K17").Value = 0.18
$ws.Range("M17").Value = 0.02
$ws.Range("O17").Value = 0.06
$ws.Range("S17").Value = 0.12
$ws.Range("H18").Value = 0.05882352941176471
$ws.Range("J18").Value = 0.5294117647058824
$ws.Range("K18").Value = 0.1176470588235294
$ws.Range("O18").Value = 0.1176470588235294
$ws.Range("S18").Value = 0.1764705882352941
$ws.Range("H19").Value = 0.2051282051282051
$ws.Range("I19").Value = 0.07692307692307693
$ws.Range("J19").Value = 0.4487179487179487
$ws.Range("K19").Value = 0.1153846153846154
$ws.Range("M19").Value = 0.02564102564102564
$ws.Range("O19").Value = 0.01282051282051282
$ws.Range("S19").Value = 0.1153846153846154
